$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 100002504
$ws.Range("J17").Value = 2782.5557
$ws.Range("L17").Value = 8347.667099999999
$ws.Range("N17").Value = -8683.667099999999
$ws.Range("H18").Value = 200
$ws.Range("I18").Value = 200
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 200
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 84
$ws.Range("N18").ClearContents()
$ws.Range("H53").Value = 2096
$ws.Range("I53").Value = 1150
$ws.Range("K53").Value = 1150
$ws.Range("M53").Value = -513
$ws.Range("H64").Value = 7455.952
$ws.Range("J64").Value = 9222.223
$ws.Range("L64").Value = 9222.223
$ws.Range("N64").Value = -9718.223
$ws.Range("H67").Value = 7455.952
$ws.Range("J67").Value = 9222.223
$ws.Range("L67").Value = 9222.223
$ws.Range("N67").Value = -10938.223
$ws.Range("H76").Value = 7535.875
$ws.Range("I76").Value = 6898.143
$ws.Range("K76").Value = 6898.143
$ws.Range("M76").Value = -6583.143
$ws.Range("H79").Value = 7535.875
$ws.Range("I79").Value = 6898.143
$ws.Range("K79").Value = 6898.143
$ws.Range("M79").Value = -5806.143
$ws.Range("H86").Value = 2994.3333
$ws.Range("I86").Value = 3390
$ws.Range("K86").Value = 3390
$ws.Range("M86").Value = -2267
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H88").Value = 3647.25
$ws.Range("I88").Value = 10194.5
$ws.Range("J88").Value = 1464.8334
$ws.Range("K88").Value = 10194.5
$ws.Range("L88").Value = 1464.8334
$ws.Range("M88").Value = -9788.5
$ws.Range("N88").Value = -2276.8334
$ws.Range("H89").Value = 2994.3333
$ws.Range("I89").Value = 3390
$ws.Range("K89").Value = 16950
$ws.Range("M89").Value = -11334
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H91").Value = 3647.25
$ws.Range("I91").Value = 10194.5
$ws.Range("J91").Value = 1464.8334
$ws.Range("K91").Value = 10194.5
$ws.Range("L91").Value = 1464.8334
$ws.Range("M91").Value = -8790.5
$ws.Range("N91").Value = -4272.8334
$ws.Range("H96").Value = 295.4375
$ws.Range("I96").Value = 235.4
$ws.Range("J96").Value = 395.5
$ws.Range("K96").Value = 706.2
$ws.Range("L96").Value = 1186.5
$ws.Range("M96").Value = 666.8
$ws.Range("N96").Value = -3932.5
$ws.Range("H100").Value = 5187.304
$ws.Range("I100").Value = 3044.1428
$ws.Range("J100").Value = 6124.9375
$ws.Range("K100").Value = 3044.1428
$ws.Range("L100").Value = 6124.9375
$ws.Range("M100").Value = -2503.1428
$ws.Range("N100").Value = -7206.9375
$ws.Range("H101").Value = 2733.25
$ws.Range("I101").Value = 1403.2
$ws.Range("J101").Value = 4950
$ws.Range("K101").Value = 4209.6
$ws.Range("L101").Value = 14850
$ws.Range("M101").Value = -2587.6
$ws.Range("N101").Value = -18094
$ws.Range("H103").Value = 1259.8
$ws.Range("I103").Value = 400
$ws.Range("J103").Value = 1474.75
$ws.Range("K103").Value = 1200
$ws.Range("L103").Value = 4424.25
$ws.Range("M103").Value = -614
$ws.Range("N103").Value = -5596.25
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H115").Value = 1797.4
$ws.Range("I115").Value = 1830.6666
$ws.Range("J115").Value = 1747.5
$ws.Range("K115").Value = 5491.9998
$ws.Range("L115").Value = 5242.5
$ws.Range("M115").Value = -3924.9998
$ws.Range("N115").Value = -8376.5
$ws.Range("H123").Value = 59999.09
$ws.Range("J123").Value = 59999.09
$ws.Range("L123").Value = 59999.09
$ws.Range("N123").Value = -69799.09
$ws.Range("H125").Value = 2201.4119
$ws.Range("I125").Value = 938.55554
$ws.Range("J125").Value = 3622.125
$ws.Range("K125").Value = 8446.99986
$ws.Range("L125").Value = 32599.125
$ws.Range("M125").Value = -5986.99986
$ws.Range("N125").Value = -37519.125
$ws.Range("H126").Value = 70000.164
$ws.Range("J126").Value = 70000.164
$ws.Range("L126").Value = 70000.164
$ws.Range("N126").Value = -79880.164
$ws.Range("H132").Value = 960.54
$ws.Range("J132").Value = 959.125
$ws.Range("L132").Value = 2877.375
$ws.Range("N132").Value = -7937.375
$ws.Range("H135").Value = 1393.44
$ws.Range("I135").Value = 1210.3043
$ws.Range("K135").Value = 10892.7387
$ws.Range("M135").Value = -8357.7387
$ws.Range("H137").Value = 3452.1724
$ws.Range("I137").Value = 3098.9092
$ws.Range("K137").Value = 9296.7276
$ws.Range("M137").Value = -6746.7276
$ws.Range("H138").Value = 2493.5745
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 2493.5745
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 7480.7235
$ws.Range("N138").Value = -17760.7235
$ws.Range("M138").ClearContents()
$ws.Range("H141").Value = 3802.1155
$ws.Range("I141").Value = 3923.2
$ws.Range("J141").Value = 775
$ws.Range("K141").Value = 11769.6
$ws.Range("L141").Value = 2325
$ws.Range("M141").Value = -6589.599999999999
$ws.Range("N141").Value = -12685

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13358.884
$ws.Range("I32").Value = 8445.556
$ws.Range("J32").Value = 41786
$ws.Range("K32").Value = 8445.556
$ws.Range("L32").Value = 41786
$ws.Range("M32").Value = -8158.556
$ws.Range("N32").Value = -42360
$ws.Range("H45").Value = 5246.727
$ws.Range("I45").Value = 3785.6667
$ws.Range("J45").Value = 7000
$ws.Range("K45").Value = 3785.6667
$ws.Range("L45").Value = 7000
$ws.Range("M45").Value = -3408.6667
$ws.Range("N45").Value = -7754
$ws.Range("H61").Value = 4577
$ws.Range("I61").Value = 2654.5833
$ws.Range("K61").Value = 2654.5833
$ws.Range("M61").Value = -2442.5833
$ws.Range("H88").Value = 10444380
$ws.Range("I88").Value = 23998.8
$ws.Range("K88").Value = 23998.8
$ws.Range("M88").Value = -23592.8
$ws.Range("H91").Value = 10444380
$ws.Range("I91").Value = 23998.8
$ws.Range("K91").Value = 23998.8
$ws.Range("M91").Value = -22594.8
$ws.Range("H97").Value = 787.59375
$ws.Range("I97").Value = 774.2
$ws.Range("K97").Value = 774.2
$ws.Range("M97").Value = -278.2
$ws.Range("H102").Value = 13890036
$ws.Range("I102").Value = 871.5
$ws.Range("K102").Value = 871.5
$ws.Range("M102").Value = 750.5
$ws.Range("H110").Value = 2195.087
$ws.Range("I110").Value = 2304.4
$ws.Range("J110").Value = 1466.3334
$ws.Range("K110").Value = 2304.4
$ws.Range("L110").Value = 1466.3334
$ws.Range("M110").Value = -259.4000000000001
$ws.Range("N110").Value = -5556.3334
$ws.Range("H122").Value = 3390.3704
$ws.Range("I122").Value = 2406.0476
$ws.Range("J122").Value = 6835.5
$ws.Range("K122").Value = 7218.1428
$ws.Range("L122").Value = 20506.5
$ws.Range("M122").Value = -4768.1428
$ws.Range("N122").Value = -25406.5
$ws.Range("H132").Value = 4381.68
$ws.Range("I132").Value = 3920.6943
$ws.Range("J132").Value = 5567.0713
$ws.Range("K132").Value = 11762.0829
$ws.Range("L132").Value = 16701.2139
$ws.Range("M132").Value = -9232.082900000001
$ws.Range("N132").Value = -21761.2139
$ws.Range("H136").Value = 4577
$ws.Range("I136").Value = 2654.5833
$ws.Range("K136").Value = 7963.749899999999
$ws.Range("M136").Value = -5413.749899999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3187.5
$ws.Range("I86").Value = 3341.6667
$ws.Range("K86").Value = 3341.6667
$ws.Range("M86").Value = -2218.6667
$ws.Range("H89").Value = 3187.5
$ws.Range("I89").Value = 3341.6667
$ws.Range("K89").Value = 16708.3335
$ws.Range("M89").Value = -11092.3335
$ws.Range("H94").Value = 7144671
$ws.Range("I94").Value = 1332.8572
$ws.Range("K94").Value = 1332.8572
$ws.Range("M94").Value = -881.8571999999999
$ws.Range("H99").Value = 1130
$ws.Range("I99").Value = 1130
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1130
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 368
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 927.26666
$ws.Range("I107").Value = 828.1818
$ws.Range("J107").Value = 1199.75
$ws.Range("K107").Value = 828.1818
$ws.Range("L107").Value = 1199.75
$ws.Range("M107").Value = 1091.8182
$ws.Range("N107").Value = -5039.75
$ws.Range("H134").Value = 3956.6296
$ws.Range("I134").Value = 3158.95
$ws.Range("J134").Value = 6235.7144
$ws.Range("K134").Value = 9476.849999999999
$ws.Range("L134").Value = 18707.1432
$ws.Range("M134").Value = -6941.849999999999
$ws.Range("N134").Value = -23777.1432

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4006.3872
$ws.Range("I31").Value = 2836.6875
$ws.Range("K31").Value = 2836.6875
$ws.Range("M31").Value = -2541.6875
$ws.Range("H34").Value = 4006.3872
$ws.Range("I34").Value = 2836.6875
$ws.Range("K34").Value = 2836.6875
$ws.Range("M34").Value = -2634.6875
$ws.Range("H99").Value = 8980310
$ws.Range("I99").Value = 1630488.4
$ws.Range("J99").Value = 20005042
$ws.Range("K99").Value = 1630488.4
$ws.Range("L99").Value = 20005042
$ws.Range("M99").Value = -1628990.4
$ws.Range("N99").Value = -20008038
$ws.Range("H126").Value = 8980310
$ws.Range("I126").Value = 1630488.4
$ws.Range("J126").Value = 20005042
$ws.Range("K126").Value = 4891465.199999999
$ws.Range("L126").Value = 60015126
$ws.Range("M126").Value = -4888995.199999999
$ws.Range("N126").Value = -60020066
$ws.Range("H132").Value = 2753.027
$ws.Range("I132").Value = 2155.5806
$ws.Range("K132").Value = 6466.7418
$ws.Range("M132").Value = -3936.7418
$ws.Range("H134").Value = 3806.2903
$ws.Range("I134").Value = 2211.7144
$ws.Range("K134").Value = 6635.1432
$ws.Range("M134").Value = -4100.1432

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 242.1
$ws.Range("I2").Value = 199.5
$ws.Range("J2").Value = 252.75
$ws.Range("K2").Value = 1197
$ws.Range("L2").Value = 1516.5
$ws.Range("M2").Value = -1084
$ws.Range("N2").Value = -1742.5
$ws.Range("H17").Value = 8674.5
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H81").Value = 2013.8334
$ws.Range("I81").Value = 833.6
$ws.Range("J81").Value = 2856.8572
$ws.Range("K81").Value = 2500.8
$ws.Range("L81").Value = 8570.571599999999
$ws.Range("M81").Value = -1377.8
$ws.Range("N81").Value = -10816.5716
$ws.Range("H84").Value = 2013.8334
$ws.Range("I84").Value = 833.6
$ws.Range("J84").Value = 2856.8572
$ws.Range("K84").Value = 7502.400000000001
$ws.Range("L84").Value = 25711.7148
$ws.Range("M84").Value = -1886.400000000001
$ws.Range("N84").Value = -36943.7148
$ws.Range("H113").Value = 1572.4445
$ws.Range("I113").Value = 1703.5
$ws.Range("K113").Value = 5110.5
$ws.Range("M113").Value = -2940.5
$ws.Range("H122").Value = 1120.7667
$ws.Range("I122").Value = 1059.5
$ws.Range("K122").Value = 9535.5
$ws.Range("M122").Value = -7085.5
$ws.Range("H132").Value = 3435.4614
$ws.Range("J132").Value = 4355.222
$ws.Range("L132").Value = 39196.998
$ws.Range("N132").Value = -44256.998
$ws.Range("H137").Value = 1748.6666
$ws.Range("J137").Value = 1999
$ws.Range("L137").Value = 5997
$ws.Range("N137").Value = -16197

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10791323
$ws.Range("I80").Value = 53287.137
$ws.Range("J80").Value = 37039856
$ws.Range("K80").Value = 53287.137
$ws.Range("L80").Value = 37039856
$ws.Range("M80").Value = -52289.137
$ws.Range("N80").Value = -37041852
$ws.Range("H83").Value = 10791323
$ws.Range("I83").Value = 53287.137
$ws.Range("J83").Value = 37039856
$ws.Range("K83").Value = 266435.685
$ws.Range("L83").Value = 185199280
$ws.Range("M83").Value = -261443.685
$ws.Range("N83").Value = -185209264
$ws.Range("H102").Value = 2703.2368
$ws.Range("I102").Value = 2109.606
$ws.Range("K102").Value = 2109.606
$ws.Range("M102").Value = -487.6060000000002
$ws.Range("H107").Value = 243.23529
$ws.Range("I107").Value = 258.46155
$ws.Range("J107").Value = 193.75
$ws.Range("K107").Value = 258.46155
$ws.Range("L107").Value = 193.75
$ws.Range("M107").Value = 1661.53845
$ws.Range("N107").Value = -4033.75
$ws.Range("H122").Value = 3987.8
$ws.Range("J122").Value = 8000
$ws.Range("L122").Value = 24000
$ws.Range("N122").Value = -28900
$ws.Range("H126").Value = 3376.7334
$ws.Range("I126").Value = 1554.5714
$ws.Range("J126").Value = 4971.125
$ws.Range("K126").Value = 4663.7142
$ws.Range("L126").Value = 14913.375
$ws.Range("M126").Value = -2193.7142
$ws.Range("N126").Value = -19853.375
$ws.Range("H132").Value = 4502.9644
$ws.Range("I132").Value = 4603.4443
$ws.Range("J132").Value = 1790
$ws.Range("K132").Value = 13810.3329
$ws.Range("L132").Value = 5370
$ws.Range("M132").Value = -11280.3329
$ws.Range("N132").Value = -10430

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3606.75
$ws.Range("I7").Value = 1994.5333
$ws.Range("J7").Value = 5029.294
$ws.Range("K7").Value = 1994.5333
$ws.Range("L7").Value = 5029.294
$ws.Range("M7").Value = -1882.5333
$ws.Range("N7").Value = -5253.294
$ws.Range("H22").Value = 1037.5714
$ws.Range("I22").Value = 1115.5
$ws.Range("J22").Value = 933.6667
$ws.Range("K22").Value = 1115.5
$ws.Range("L22").Value = 933.6667
$ws.Range("M22").Value = -820.5
$ws.Range("N22").Value = -1523.6667
$ws.Range("H27").Value = 1037.5714
$ws.Range("I27").Value = 1115.5
$ws.Range("J27").Value = 933.6667
$ws.Range("K27").Value = 1115.5
$ws.Range("L27").Value = 933.6667
$ws.Range("M27").Value = -1008.5
$ws.Range("N27").Value = -1147.6667
$ws.Range("H40").Value = 6038.3945
$ws.Range("I40").Value = 5937.5356
$ws.Range("J40").Value = 6320.8
$ws.Range("K40").Value = 5937.5356
$ws.Range("L40").Value = 6320.8
$ws.Range("M40").Value = -5801.5356
$ws.Range("N40").Value = -6592.8
$ws.Range("H80").Value = 108620
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 108620
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 108620
$ws.Range("N80").Value = -110866
$ws.Range("M80").ClearContents()
$ws.Range("H82").Value = 1927.6364
$ws.Range("I82").Value = 1149.25
$ws.Range("J82").Value = 2372.4285
$ws.Range("K82").Value = 1149.25
$ws.Range("L82").Value = 2372.4285
$ws.Range("M82").Value = -788.25
$ws.Range("N82").Value = -3094.4285
$ws.Range("H83").Value = 108620
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 108620
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 325860
$ws.Range("N83").Value = -337092
$ws.Range("M83").ClearContents()
$ws.Range("H85").Value = 1927.6364
$ws.Range("I85").Value = 1149.25
$ws.Range("J85").Value = 2372.4285
$ws.Range("K85").Value = 1149.25
$ws.Range("L85").Value = 2372.4285
$ws.Range("M85").Value = 98.75
$ws.Range("N85").Value = -4868.4285
$ws.Range("H126").Value = 3606.75
$ws.Range("I126").Value = 1994.5333
$ws.Range("J126").Value = 5029.294
$ws.Range("K126").Value = 5983.5999
$ws.Range("L126").Value = 15087.882
$ws.Range("M126").Value = -3513.5999
$ws.Range("N126").Value = -20027.882
$ws.Range("H132").Value = 4493
$ws.Range("I132").Value = 3683.1853
$ws.Range("K132").Value = 11049.5559
$ws.Range("M132").Value = -8519.555899999999
$ws.Range("H136").Value = 4693.9653
$ws.Range("I136").Value = 3298.75
$ws.Range("K136").Value = 9896.25
$ws.Range("M136").Value = -7346.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 30495
$ws.Range("J40").Value = 30495
$ws.Range("L40").Value = 30495
$ws.Range("N40").Value = -30793
$ws.Range("H74").Value = 12560
$ws.Range("J74").Value = 13415
$ws.Range("L74").Value = 13415
$ws.Range("N74").Value = -15287
$ws.Range("H77").Value = 12560
$ws.Range("J77").Value = 13415
$ws.Range("L77").Value = 40245
$ws.Range("N77").Value = -49605
$ws.Range("H81").Value = 18522040
$ws.Range("I81").Value = 2360.75
$ws.Range("K81").Value = 4721.5
$ws.Range("M81").Value = -3660.5
$ws.Range("H84").Value = 18522040
$ws.Range("I84").Value = 2360.75
$ws.Range("K84").Value = 23607.5
$ws.Range("M84").Value = -18303.5
$ws.Range("H96").Value = 31523.766
$ws.Range("I96").Value = 40424.46
$ws.Range("J96").Value = 2596.5
$ws.Range("K96").Value = 40424.46
$ws.Range("L96").Value = 2596.5
$ws.Range("M96").Value = -39051.46
$ws.Range("N96").Value = -5342.5
$ws.Range("H113").Value = 557.1111
$ws.Range("J113").Value = 485.33334
$ws.Range("L113").Value = 1456.00002
$ws.Range("N113").Value = -5796.000019999999
$ws.Range("H132").Value = 1610.6459
$ws.Range("I132").Value = 1039.6744
$ws.Range("J132").Value = 6521
$ws.Range("K132").Value = 3119.023200000001
$ws.Range("L132").Value = 19563
$ws.Range("M132").Value = -589.0232000000005
$ws.Range("N132").Value = -24623
$ws.Range("H136").Value = 10873410
$ws.Range("I136").Value = 18521606
$ws.Range("J136").Value = 4922.263
$ws.Range("K136").Value = 55564818
$ws.Range("L136").Value = 14766.789
$ws.Range("M136").Value = -55562268
$ws.Range("N136").Value = -19866.789

Write-Host "Applied all changes"